# Regenerate the "K" column (column G) values on Sheet1.
# The data-generation script was re-run ("regen save_data to use K instead
# of Strike#, regen std/mean, calc and write s_vals"), which produced new
# strikeout ("K") counts for every existing row (rows 2-70). Only column G
# values change; everything else on the sheet stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2 through 70, in order.
$newK = @(2,1,0,1,1,1,0,1,2,0,1,1,2,0,1,2,2,3,3,0,1,0,0,0,1,2,5,2,0,1,2,4,2,1,3,1,1,0,0,2,1,1,0,1,1,0,0,1,0,2,0,0,1,1,0,1,1,1,1,0,2,1,3,0,2,2,1,4,1)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
